# Use conventional styles/indents for Word bullet lists.
#
# The document's numbering definition (word/numbering.xml) defines a
# 9-level bullet abstract numbering where every level uses a 480-twip
# hanging indent. Word's "conventional" bullet list indents use a
# 360-twip hanging indent instead, so update w:hanging from 480 to 360
# on every w:ind element in that numbering definition (the w:left
# values themselves are unchanged).
#
# There is no dedicated ListLevels/ListTemplate setter that reaches the
# underlying numbering.xml part in this runtime, so we round-trip the
# package through Document.WordOpenXML (the flat-OPC representation of
# the whole .docx, including word/numbering.xml) and patch the
# w:hanging="480" occurrences there.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML
$xml = $xml.Replace('w:hanging="480"', 'w:hanging="360"')
$d.WordOpenXML = $xml

Write-Output "updated hanging indents"
